$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the date
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2022-11-08", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2022-11-09", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. First new section: an "Introduction" body paragraph right after the
#    first (empty) Heading 1 paragraph that currently sits inside the
#    "section" bookmark.
# ---------------------------------------------------------------------------
$pIntroHeading = $d.Paragraphs.Item(4)
$pIntroHeading.Range.InsertParagraphAfter()

$pIntro = $d.Paragraphs.Item(5)
$pIntro.Style = "FirstParagraph"
$pIntro.Range.Text = "Predators and prey have a suite of strategies they use to maximize " + `
    "benefit with the least cost. Predators are able to switch between food " + `
    "sources when one source becomes hard to find (Carle and Rowe 2014). Other " + `
    "strategies can be gleaned from our foraging lab - central place foraging " + `
    "and traplining, in which a predator will center themselves where prey " + `
    "generally congregates or when a predator uses the most efficient path " + `
    "possible. Prey strategies focus decreasing the cost of predation by " + `
    "predator avoidance, group foraging, or mimicry. In these situations prey " + `
    "can avoid capture, work as a group, or mimic a toxic species. Each " + `
    "strategy, either predator or prey, confers selective pressures on the " + `
    "other groups and on associated species not involved in the interaction."

# ---------------------------------------------------------------------------
# 3. The paragraph that used to read "References" becomes the second empty
#    Heading 1 (its own new body paragraph follows), since the References
#    section is being rebuilt further down.
# ---------------------------------------------------------------------------
$pOldReferences = $d.Paragraphs.Item(6)
$clearRange = $d.Range($pOldReferences.Range.Start, $pOldReferences.Range.End - 1)
$clearRange.Delete()

$pOldReferences.Range.InsertParagraphAfter()
$pSection1Body = $d.Paragraphs.Item(7)
$pSection1Body.Style = "FirstParagraph"
$pSection1Body.Range.Text = "There is evidence that predators will switch their food source to a " + `
    "toxic prey when undefended prey are better concealed (Carle and Rowe " + `
    "2014), however the authors stipulate a natural environment could " + `
    "undermine this switch in diet. Their results suggest that this change " + `
    "in diet follows a risk-prone strategy, in which, the predator eats more " + `
    "toxic prey when the probability of finding undefended prey is low. This " + `
    "strategy leads to selective pressure on the predator to tolerate toxins " + `
    "and on the prey to be more cryptic or more toxic. In the foraging lab, " + `
    "simulated predators followed two strategies - central place foraging " + `
    "and traplining - in the first, one predator stayed where prey was most " + `
    "abundant while the latter predator took the most efficient path to " + `
    "capture prey. For central place foraging, this strategy maximizes " + `
    "benefit and minimizes the cost of searching for food. Alternatively, " + `
    "traplining does the same to a lesser degree by maximizing prey capture " + `
    "and reduces the cost to search."

# ---------------------------------------------------------------------------
# 4. New "Prey Strategies" and "Conclusion" headings.
# ---------------------------------------------------------------------------
$pSection1Body.Range.InsertParagraphAfter()
$pPreyStrategies = $d.Paragraphs.Item(8)
$pPreyStrategies.Style = "Heading1"
$pPreyStrategies.Range.Text = "Prey Strategies"

$pPreyStrategies.Range.InsertParagraphAfter()
$pConclusion = $d.Paragraphs.Item(9)
$pConclusion.Style = "Heading1"
$pConclusion.Range.Text = "Conclusion"

# ---------------------------------------------------------------------------
# 5. New "References" heading plus the bibliography entry (with hyperlink).
# ---------------------------------------------------------------------------
$pConclusion.Range.InsertParagraphAfter()
$pReferences = $d.Paragraphs.Item(10)
$pReferences.Style = "Heading1"
$pReferences.Range.Text = "References"

$pReferences.Range.InsertParagraphAfter()
$pBib = $d.Paragraphs.Item(11)
$pBib.Style = "Bibliography"
$enDash = [char]0x2013
$pBib.Range.Text = "Carle, T., and C. Rowe. 2014. Avian predators change their foraging " + `
    "strategy on defended prey when undefended prey are hard to find. " + `
    "Animal Behaviour 93:97" + $enDash + "103."

$bibRange = $d.Range($pBib.Range.Start, $pBib.Range.End)
$bibRange.Find.Execute("Avian predators change their foraging strategy on defended prey when undefended prey are hard to find", `
                        $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($bibRange, "https://doi.org/10.1016/j.anbehav.2014.04.025") | Out-Null

# ---------------------------------------------------------------------------
# 6. Bookmarks (names matter; the numeric ids Word assigns on save are an
#    internal implementation detail we do not control from the object
#    model, so we only manage the names/spans here).
# ---------------------------------------------------------------------------
$introHeading = $d.Paragraphs.Item(4)
$introBody = $d.Paragraphs.Item(5)
$sectionRange = $d.Range($introHeading.Range.Start, $introBody.Range.End)
$d.Bookmarks.Add("section", $sectionRange) | Out-Null

$section1Heading = $d.Paragraphs.Item(6)
$section1Body = $d.Paragraphs.Item(7)
$section1Range = $d.Range($section1Heading.Range.Start, $section1Body.Range.End)
$d.Bookmarks.Add("section-1", $section1Range) | Out-Null

$preyStrategiesPara = $d.Paragraphs.Item(8)
$preyRange = $d.Range($preyStrategiesPara.Range.Start, $preyStrategiesPara.Range.End)
$d.Bookmarks.Add("prey-strategies", $preyRange) | Out-Null

$conclusionPara = $d.Paragraphs.Item(9)
$conclusionRange = $d.Range($conclusionPara.Range.Start, $conclusionPara.Range.End)
$d.Bookmarks.Add("conclusion", $conclusionRange) | Out-Null

$referencesHeadingPara = $d.Paragraphs.Item(10)
$bibPara = $d.Paragraphs.Item(11)
$referencesRange = $d.Range($referencesHeadingPara.Range.Start, $bibPara.Range.End)
$d.Bookmarks.Add("references", $referencesRange) | Out-Null

$refsRange = $d.Range($bibPara.Range.Start, $bibPara.Range.End)
$d.Bookmarks.Add("refs", $refsRange) | Out-Null
$d.Bookmarks.Add("ref-carle", $refsRange) | Out-Null

Write-Output "edit complete"
